$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily auto-update job appends the newest draw result as a new row at the
# bottom of the "Results" sheet.
$row = 76

# Columns A (date) and C (6-digit phase code) look numeric/date-like, so
# Excel would otherwise auto-convert them on assignment. Force text
# formatting first so they stay literal text, then restore the default
# "Normal" style afterwards so the cell doesn't end up tagged with a
# lingering custom number format (matching the source feed, which stores
# every column as plain text with no special formatting).
$ws.Range("A" + $row).NumberFormat = "@"
$ws.Range("C" + $row).NumberFormat = "@"

$ws.Range("A" + $row).Value = "2025-12-01"
$ws.Range("B" + $row).Value = "Pick 3"
$ws.Range("C" + $row).Value = "251201"
$ws.Range("D" + $row).Value = "4-7-4"
$ws.Range("E" + $row).Value = "2025-12-01T21:46:49.626+04:00"

$ws.Range("A" + $row).Style = "Normal"
$ws.Range("C" + $row).Style = "Normal"
